$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$data = @(
  @{C='1050960713'; D='SOLBEY MILENA BOLAÑO MORALES'; E='2403'; F=40533; G=1600000},
  @{C='1050960713'; D='SOLBEY MILENA BOLAÑO MORALES'; E='2402'; F=64000; G=1600000},
  @{C='1050960713'; D='SOLBEY MILENA BOLAÑO MORALES'; E='2401'; F=64000; G=1600000},
  @{C='1050960713'; D='SOLBEY MILENA BOLAÑO MORALES'; E='2312'; F=64000; G=1600000},
  @{C='1050960713'; D='SOLBEY MILENA BOLAÑO MORALES'; E='2311'; F=64000; G=1600000},
  @{C='1050960713'; D='SOLBEY MILENA BOLAÑO MORALES'; E='2310'; F=64000; G=1600000},
  @{C='1050960713'; D='SOLBEY MILENA BOLAÑO MORALES'; E='2309'; F=64000; G=1600000},
  @{C='1026298364'; D='MELANIE SILVANA SANDRINO RIOS'; E='2403'; F=45600; G=1800000},
  @{C='1026298364'; D='MELANIE SILVANA SANDRINO RIOS'; E='2402'; F=72000; G=1800000},
  @{C='1026298364'; D='MELANIE SILVANA SANDRINO RIOS'; E='2401'; F=72000; G=1800000},
  @{C='1026298364'; D='MELANIE SILVANA SANDRINO RIOS'; E='2312'; F=72000; G=1800000},
  @{C='1026298364'; D='MELANIE SILVANA SANDRINO RIOS'; E='2311'; F=72000; G=1800000},
  @{C='1026298364'; D='MELANIE SILVANA SANDRINO RIOS'; E='2310'; F=72000; G=1800000},
  @{C='1026298364'; D='MELANIE SILVANA SANDRINO RIOS'; E='2309'; F=72000; G=1800000},
  @{C='1026298364'; D='MELANIE SILVANA SANDRINO RIOS'; E='2308'; F=72000; G=1800000},
  @{C='1026298364'; D='MELANIE SILVANA SANDRINO RIOS'; E='2307'; F=72000; G=1800000},
  @{C='1026298364'; D='MELANIE SILVANA SANDRINO RIOS'; E='2306'; F=72000; G=1800000},
  @{C='1026298364'; D='MELANIE SILVANA SANDRINO RIOS'; E='2305'; F=72000; G=1800000},
  @{C='1026298364'; D='MELANIE SILVANA SANDRINO RIOS'; E='2304'; F=72000; G=1800000},
  @{C='1026298364'; D='MELANIE SILVANA SANDRINO RIOS'; E='2303'; F=72000; G=1800000},
  @{C='1026298364'; D='MELANIE SILVANA SANDRINO RIOS'; E='2302'; F=72000; G=1800000},
  @{C='1026298364'; D='MELANIE SILVANA SANDRINO RIOS'; E='2301'; F=72000; G=1800000},
  @{C='1026298364'; D='MELANIE SILVANA SANDRINO RIOS'; E='2212'; F=72000; G=1800000},
  @{C='1026298364'; D='MELANIE SILVANA SANDRINO RIOS'; E='2211'; F=72000; G=1800000},
  @{C='1026298364'; D='MELANIE SILVANA SANDRINO RIOS'; E='2210'; F=64800; G=1800000},
  @{C='1096253977'; D='SARETH CASTRO MONTES'; E='2403'; F=40533; G=1600000},
  @{C='1096253977'; D='SARETH CASTRO MONTES'; E='2402'; F=64000; G=1600000},
  @{C='1096253977'; D='SARETH CASTRO MONTES'; E='2401'; F=64000; G=1600000},
  @{C='1096253977'; D='SARETH CASTRO MONTES'; E='2312'; F=64000; G=1600000},
  @{C='1096253977'; D='SARETH CASTRO MONTES'; E='2311'; F=64000; G=1600000},
  @{C='1096253977'; D='SARETH CASTRO MONTES'; E='2310'; F=64000; G=1600000},
  @{C='1096253977'; D='SARETH CASTRO MONTES'; E='2309'; F=64000; G=1600000},
  @{C='1096253977'; D='SARETH CASTRO MONTES'; E='2308'; F=64000; G=1600000},
  @{C='1096253977'; D='SARETH CASTRO MONTES'; E='2307'; F=64000; G=1600000},
  @{C='1096253977'; D='SARETH CASTRO MONTES'; E='2306'; F=64000; G=1600000},
  @{C='1096253977'; D='SARETH CASTRO MONTES'; E='2305'; F=64000; G=1600000},
  @{C='1096253977'; D='SARETH CASTRO MONTES'; E='2304'; F=64000; G=1600000},
  @{C='1096253977'; D='SARETH CASTRO MONTES'; E='2303'; F=64000; G=1600000},
  @{C='1096253977'; D='SARETH CASTRO MONTES'; E='2302'; F=64000; G=1600000},
  @{C='1096253977'; D='SARETH CASTRO MONTES'; E='2301'; F=64000; G=1600000},
  @{C='1096253977'; D='SARETH CASTRO MONTES'; E='2212'; F=64000; G=1600000},
  @{C='1096253977'; D='SARETH CASTRO MONTES'; E='2211'; F=64000; G=1600000},
  @{C='1096253977'; D='SARETH CASTRO MONTES'; E='2210'; F=57600; G=1600000}
)

$row = 16
foreach ($item in $data) {
  $ws.Cells.Item($row, 3).Formula = $item.C
  $ws.Cells.Item($row, 4).Formula = $item.D
  $ws.Cells.Item($row, 5).Formula = $item.E
  $ws.Cells.Item($row, 6).Formula = $item.F
  $ws.Cells.Item($row, 7).Formula = $item.G
  $row = $row + 1
}

Write-Host "Updated rows 16 to $($row - 1)"
